# Refresh the crypto symbol list (coinranking.com snapshot) with the
# prices, 24h-volume labels, and "Hora" (hour) column captured by the
# "Updated symbol list" GitHub Actions run on Thu Dec 29 21:07:38 UTC 2022.
#
# Every cell on the sheet is stored as a literal/inline text string (the
# sheet uses t="inlineStr" throughout -- even for numeric-looking values
# such as prices or the hour counter in column G). A plain
# `Range.Value = "..."` assignment is not enough for those cells: like
# real Excel, the COM layer auto-coerces a numeric-looking string into a
# genuine number (and the xlsx would then serialize it as t="n"). To keep
# the original text semantics we briefly force the cell to Text format,
# assign the literal string, then restore the cell's style so no
# formatting changes leak into the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "245.74"; NumericLooking = $true },
    @{ Cell = "G2"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D3"; Value = "24.17"; NumericLooking = $true },
    @{ Cell = "G3"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G4"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D5"; Value = "0.05767"; NumericLooking = $true },
    @{ Cell = "G5"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D6"; Value = "6.461"; NumericLooking = $true },
    @{ Cell = "G6"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D7"; Value = "3.148"; NumericLooking = $true },
    @{ Cell = "G7"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D8"; Value = "0.8166"; NumericLooking = $true },
    @{ Cell = "G8"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D9"; Value = "0.8469"; NumericLooking = $true },
    @{ Cell = "G9"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D10"; Value = "0.1359"; NumericLooking = $true },
    @{ Cell = "G10"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D11"; Value = "0.06943"; NumericLooking = $true },
    @{ Cell = "G11"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D12"; Value = "0.03135"; NumericLooking = $true },
    @{ Cell = "G12"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D13"; Value = "0.02927"; NumericLooking = $true },
    @{ Cell = "G13"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D14"; Value = "0.09391"; NumericLooking = $true },
    @{ Cell = "G14"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D15"; Value = "3.750"; NumericLooking = $true },
    @{ Cell = "G15"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G16"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G17"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D18"; Value = "0.0005983"; NumericLooking = $true },
    @{ Cell = "G18"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D19"; Value = "0.006190"; NumericLooking = $true },
    @{ Cell = "G19"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D20"; Value = "0.001236"; NumericLooking = $true },
    @{ Cell = "G20"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D21"; Value = "0.004618"; NumericLooking = $true },
    @{ Cell = "E21"; Value = "20HotbitTokenHTBBestin24h"; NumericLooking = $false },
    @{ Cell = "G21"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D22"; Value = "0.00006900"; NumericLooking = $true },
    @{ Cell = "G22"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D23"; Value = "3.499"; NumericLooking = $true },
    @{ Cell = "G23"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G24"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D25"; Value = "0.3196"; NumericLooking = $true },
    @{ Cell = "G25"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D26"; Value = "0.1302"; NumericLooking = $true },
    @{ Cell = "G26"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D27"; Value = "0.1360"; NumericLooking = $true },
    @{ Cell = "G27"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D28"; Value = "0.0002333"; NumericLooking = $true },
    @{ Cell = "G28"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G29"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G30"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G31"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G32"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G33"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G34"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G35"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G36"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G37"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G38"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G39"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D40"; Value = "0.03654"; NumericLooking = $true },
    @{ Cell = "G40"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D41"; Value = "0.006270"; NumericLooking = $true },
    @{ Cell = "E41"; Value = "40KickTokenKICK"; NumericLooking = $false },
    @{ Cell = "G41"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D42"; Value = "0.1055"; NumericLooking = $true },
    @{ Cell = "G42"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D43"; Value = "0.002730"; NumericLooking = $true },
    @{ Cell = "G43"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G44"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D45"; Value = "0.00005255"; NumericLooking = $true },
    @{ Cell = "G45"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G46"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D47"; Value = "0.3771"; NumericLooking = $true },
    @{ Cell = "G47"; Value = "21"; NumericLooking = $true },
    @{ Cell = "D48"; Value = "0.002291"; NumericLooking = $true },
    @{ Cell = "G48"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G49"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G50"; Value = "21"; NumericLooking = $true },
    @{ Cell = "G51"; Value = "21"; NumericLooking = $true }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.NumericLooking) {
        # Force text storage so "245.74" / "21" / "0.00006900" etc. stay
        # strings instead of being auto-converted to numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
